# Adding to goals section.
# Move/resize the two translucent "cover" rectangles ("Rectangle 1" and
# "Rectangle 30") on the goals mockup slide further along the timeline.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$rect1 = $s.Shapes.Item("Rectangle 1")
$rect1.Left   = 503.9999541598425
$rect1.Top    = 328.9886168771654
$rect1.Width  = 45.7500381
$rect1.Height = 133.01145932283464

$rect30 = $s.Shapes.Item("Rectangle 30")
$rect30.Left   = 549.7499389598424
$rect30.Top    = 283.20413208818894
$rect30.Width  = 39.000433000787396
$rect30.Height = 186.0000381
